$wb = $excel.ActiveWorkbook

# Add the new worksheet "irctc" after the last existing sheet (IssueDate)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "irctc"

# Fill in the data
$newSheet.Range("A1").Value = "city"
$newSheet.Range("B1").Value = "phoneNo"
$newSheet.Range("A2").Value = "New Delhi"
$newSheet.Range("B2").Value = ": 011-23221147"

# Bold header style for row 1
$newSheet.Range("A1:B1").Font.Bold = $true

# Select B2 and activate this sheet
$newSheet.Range("B2").Select()
$newSheet.Activate()
